$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the number formatting used by the existing data rows (column A = bold
# bordered "index" style, column E = datetime style) onto the two new rows
# so the freshly-appended cells line up with the rest of the table.
$ws.Range("A2").Copy()
$ws.Range("A71").PasteSpecial(-4122)
$ws.Range("A72").PasteSpecial(-4122)

$ws.Range("E2").Copy()
$ws.Range("E71").PasteSpecial(-4122)
$ws.Range("E72").PasteSpecial(-4122)

# Row 71 - Port MTI FC vs Bangkok Utd
$ws.Cells.Item(71, 1).Value = 70
$ws.Cells.Item(71, 2).Value = "thailand"
$ws.Cells.Item(71, 3).Value = "thai-league-1"
$ws.Cells.Item(71, 4).Value = "2023-2024"
$ws.Cells.Item(71, 5).Value = 45234.5
$ws.Cells.Item(71, 6).Value = "Port MTI FC"
$ws.Cells.Item(71, 7).Value = 0
$ws.Cells.Item(71, 8).Value = "Bangkok Utd"
$ws.Cells.Item(71, 9).Value = 2
$ws.Cells.Item(71, 10).Value = 2.02
$ws.Cells.Item(71, 11).Value = "29/10/2023 18:42"
$ws.Cells.Item(71, 12).Value = 2.46
$ws.Cells.Item(71, 13).Value = "04/11/2023 11:58"
$ws.Cells.Item(71, 14).Value = 3.66
$ws.Cells.Item(71, 15).Value = "29/10/2023 18:42"
$ws.Cells.Item(71, 16).Value = 3.66
$ws.Cells.Item(71, 17).Value = "04/11/2023 11:58"
$ws.Cells.Item(71, 18).Value = 3.28
$ws.Cells.Item(71, 19).Value = "29/10/2023 18:42"
$ws.Cells.Item(71, 20).Value = 2.74
$ws.Cells.Item(71, 21).Value = "04/11/2023 11:56"
$ws.Cells.Item(71, 22).Value = "https://www.betexplorer.com/football/thailand/thai-league-1/port-mti-fc-bangkok-utd/pzevxA0C/"

# Row 72 - Khonkaen Utd. vs Uthai Thani
$ws.Cells.Item(72, 1).Value = 71
$ws.Cells.Item(72, 2).Value = "thailand"
$ws.Cells.Item(72, 3).Value = "thai-league-1"
$ws.Cells.Item(72, 4).Value = "2023-2024"
$ws.Cells.Item(72, 5).Value = 45234.54166666666
$ws.Cells.Item(72, 6).Value = "Khonkaen Utd."
$ws.Cells.Item(72, 7).Value = 2
$ws.Cells.Item(72, 8).Value = "Uthai Thani"
$ws.Cells.Item(72, 9).Value = 2
$ws.Cells.Item(72, 10).Value = 3.01
$ws.Cells.Item(72, 11).Value = "03/11/2023 15:43"
$ws.Cells.Item(72, 12).Value = 2.93
$ws.Cells.Item(72, 13).Value = "04/11/2023 12:52"
$ws.Cells.Item(72, 14).Value = 3.81
$ws.Cells.Item(72, 15).Value = "03/11/2023 15:43"
$ws.Cells.Item(72, 16).Value = 4.03
$ws.Cells.Item(72, 17).Value = "04/11/2023 12:52"
$ws.Cells.Item(72, 18).Value = 2.09
$ws.Cells.Item(72, 19).Value = "03/11/2023 15:43"
$ws.Cells.Item(72, 20).Value = 2.2
$ws.Cells.Item(72, 21).Value = "04/11/2023 12:52"
$ws.Cells.Item(72, 22).Value = "https://www.betexplorer.com/football/thailand/thai-league-1/khonkaen-united-uthai-thani/xAiWwWVa/"
